# Apply the "rebuild for panorama snippets" edit:
#  1. Update PANORAMA_IP value on the 'values' sheet from 192.168.55.161 to 192.168.55.8
#  2. Remove the INTERNET_ZONE row (row 24) from the 'values' sheet entirely, which
#     shifts the rows below it up by one and causes all dependent formulas on the
#     'set commands' sheet (which reference 'values'!B25.. etc.) to be renumbered
#     automatically, along with the shared strings table being rebuilt on save.

$wb = $excel.ActiveWorkbook
$valuesWs = $wb.Worksheets.Item("values")

# 1. Update PANORAMA_IP value (row 4, column B)
$valuesWs.Range("B4").Value = "192.168.55.8"

# 2. Delete the INTERNET_ZONE row (row 24: INTERNET_ZONE / internet / untrust zone to filter out in reports)
$valuesWs.Rows.Item(24).Delete()
